# New table for copying range
# Adds a "Table3" table (Name / Date / Range / Price) below the existing
# PlanTable on the Summary sheet, and updates the Config sheet's
# documentation rows (view state + the "Target table" label that used to
# say "PlanTable") to refer to the new table.

$wb = $excel.ActiveWorkbook

# --- Summary sheet: header row + new table ---
$ws1 = $wb.Worksheets.Item("Summary")

$ws1.Range("B15").Value = "Name"
$ws1.Range("C15").Value = "Date"
$ws1.Range("D15").Value = "Range"
$ws1.Range("E15").Value = "Price"

$tableRange = $ws1.Range("B15:E16")
$table3 = $ws1.ListObjects.Add(1, $tableRange, $null, 1)
$table3.Name = "Table3"
$table3.HeaderRowRange.Font.Bold = $true

# --- Config sheet: point the "Replace table" example at the new table ---
$ws2 = $wb.Worksheets.Item("Config")
$ws2.Range("D36").Value = "Table3"

$ws2.Activate()
$ws2.Range("D38").Select()

# Leave the Summary sheet selected/active, matching the saved workbook state
$ws1.Activate()
$ws1.Range("B19").Select()
